$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated symbol list on Wed Jan 11 17:09:52 UTC 2023 with GitHub Actions
# Refresh Price (D), Volume(1h) (E), and Hora (G) columns for each crypto
# row with the latest scraped values from coinranking.com.
#
# Values are assigned with a leading apostrophe so Excel treats the
# numeric-looking strings (prices, percentages, hour) as literal TEXT -
# matching the existing text cells in this sheet - rather than silently
# converting them to the Number type. The trailing bulk Style reset clears
# the "quote prefix" formatting flag that Excel applies automatically when
# a value is entered with a leading apostrophe, so the cell formatting is
# left exactly as it was before the edit.

$ws.Range("D2").Value = "'276.56"
$ws.Range("E2").Value = "'0.55%"
$ws.Range("G2").Value = "'17"
$ws.Range("D3").Value = "'27.30"
$ws.Range("E3").Value = "'0.65%"
$ws.Range("G3").Value = "'17"
$ws.Range("E4").Value = "'1.29%"
$ws.Range("G4").Value = "'17"
$ws.Range("D5").Value = "'0.06326"
$ws.Range("E5").Value = "'0.37%"
$ws.Range("G5").Value = "'17"
$ws.Range("D6").Value = "'7.024"
$ws.Range("E6").Value = "'1.22%"
$ws.Range("G6").Value = "'17"
$ws.Range("D7").Value = "'1.374"
$ws.Range("E7").Value = "'3.61%"
$ws.Range("G7").Value = "'17"
$ws.Range("D8").Value = "'0.8881"
$ws.Range("E8").Value = "'1.25%"
$ws.Range("G8").Value = "'17"
$ws.Range("D9").Value = "'0.1509"
$ws.Range("E9").Value = "'-0.40%"
$ws.Range("G9").Value = "'17"
$ws.Range("D10").Value = "'0.05269"
$ws.Range("E10").Value = "'4.65%"
$ws.Range("G10").Value = "'17"
$ws.Range("D11").Value = "'0.07437"
$ws.Range("E11").Value = "'-0.87%"
$ws.Range("G11").Value = "'17"
$ws.Range("D12").Value = "'0.02889"
$ws.Range("E12").Value = "'-0.48%"
$ws.Range("G12").Value = "'17"
$ws.Range("D13").Value = "'0.08945"
$ws.Range("E13").Value = "'-0.73%"
$ws.Range("G13").Value = "'17"
$ws.Range("D14").Value = "'0.001572"
$ws.Range("E14").Value = "'0.62%"
$ws.Range("G14").Value = "'17"
$ws.Range("D15").Value = "'0.0006346"
$ws.Range("E15").Value = "'-0.54%"
$ws.Range("G15").Value = "'17"
$ws.Range("D16").Value = "'0.006033"
$ws.Range("E16").Value = "'0.25%"
$ws.Range("G16").Value = "'17"
$ws.Range("D17").Value = "'3.472"
$ws.Range("E17").Value = "'0.71%"
$ws.Range("G17").Value = "'17"
$ws.Range("D18").Value = "'3.297"
$ws.Range("E18").Value = "'-0.20%"
$ws.Range("G18").Value = "'17"
$ws.Range("D19").Value = "'2.233"
$ws.Range("E19").Value = "'-1.70%"
$ws.Range("G19").Value = "'17"
$ws.Range("D20").Value = "'0.3170"
$ws.Range("E20").Value = "'1.67%"
$ws.Range("G20").Value = "'17"
$ws.Range("E21").Value = "'0.60%"
$ws.Range("G21").Value = "'17"
$ws.Range("D22").Value = "'3.929"
$ws.Range("E22").Value = "'0.53%"
$ws.Range("G22").Value = "'17"
$ws.Range("G23").Value = "'17"
$ws.Range("E24").Value = "'-0.52%"
$ws.Range("G24").Value = "'17"
$ws.Range("D25").Value = "'0.001176"
$ws.Range("E25").Value = "'0.27%"
$ws.Range("G25").Value = "'17"
$ws.Range("D26").Value = "'0.004246"
$ws.Range("E26").Value = "'10.69%"
$ws.Range("G26").Value = "'17"
$ws.Range("G27").Value = "'17"
$ws.Range("D28").Value = "'0.0001180"
$ws.Range("E28").Value = "'-1.78%"
$ws.Range("G28").Value = "'17"
$ws.Range("E29").Value = "'-14.95%"
$ws.Range("G29").Value = "'17"
$ws.Range("G30").Value = "'17"
$ws.Range("G31").Value = "'17"
$ws.Range("G32").Value = "'17"
$ws.Range("G33").Value = "'17"
$ws.Range("G34").Value = "'17"
$ws.Range("G35").Value = "'17"
$ws.Range("G36").Value = "'17"
$ws.Range("G37").Value = "'17"
$ws.Range("G38").Value = "'17"
$ws.Range("G39").Value = "'17"
$ws.Range("D40").Value = "'0.03969"
$ws.Range("E40").Value = "'-3.18%"
$ws.Range("G40").Value = "'17"
$ws.Range("D41").Value = "'0.006662"
$ws.Range("E41").Value = "'-2.81%"
$ws.Range("G41").Value = "'17"
$ws.Range("D42").Value = "'0.1398"
$ws.Range("E42").Value = "'19.42%"
$ws.Range("G42").Value = "'17"
$ws.Range("D43").Value = "'0.001910"
$ws.Range("E43").Value = "'-13.68%"
$ws.Range("G43").Value = "'17"
$ws.Range("D44").Value = "'0.01171"
$ws.Range("E44").Value = "'1.74%"
$ws.Range("G44").Value = "'17"
$ws.Range("D45").Value = "'0.00005437"
$ws.Range("E45").Value = "'4.94%"
$ws.Range("G45").Value = "'17"
$ws.Range("E46").Value = "'5.01%"
$ws.Range("G46").Value = "'17"
$ws.Range("E47").Value = "'-19.66%"
$ws.Range("G47").Value = "'17"
$ws.Range("G48").Value = "'17"
$ws.Range("G49").Value = "'17"
$ws.Range("G50").Value = "'17"
$ws.Range("G51").Value = "'17"

# Clear the quote-prefix flag picked up from the apostrophe-prefixed
# assignments above so cell formatting is unchanged by this update.
$ws.Range("D2:G51").Style = "Normal"
